$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '37.755.88'
$ws.Range("E2").Value = '  -1.14%  '

$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '2.029.18'
$ws.Range("E3").Value = '  -1.61%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '227.36'
$ws.Range("E5").Value = '  -1.55%  '

$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '0.613'
$ws.Range("E6").Value = '  -0.53%  '

$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '59.31'
$ws.Range("E7").Value = '  +1.34%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("E9").Value = '  -0.99%  '

$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '0.0814'
$ws.Range("E10").Value = '  +0.50%  '

$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '14.57'
$ws.Range("E12").Value = '  -0.76%  '

$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '2.330.44'
$ws.Range("E13").Value = '  -1.54%  '

$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '21.01'
$ws.Range("E14").Value = '  +1.32%  '

$ws.Range("E15").Value = '  +0.24%  '

$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '5.19'
$ws.Range("E16").Value = '  -1.96%  '

$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '2.028.29'
$ws.Range("E17").Value = '  -1.80%  '

$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '37.678.21'
$ws.Range("E18").Value = '  -0.97%  '

$ws.Range("E19").Value = '  -2.09%  '

$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '69.72'
$ws.Range("E20").Value = '  -0.39%  '

$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '0.0₃0823'
$ws.Range("E21").Value = '  -1.26%  '

$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '224.70'
$ws.Range("E22").Value = '  -0.15%  '

$ws.Range("E23").Value = '  +0.05%  '

$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '2.39'
$ws.Range("E24").Value = '  -2.80%  '

$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '2.21'
$ws.Range("E25").Value = '  -1.76%  '

$ws.Range("E26").Value = '  -1.29%  '

$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '165.04'
$ws.Range("E27").Value = '  -0.78%  '

$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '0.128'
$ws.Range("E28").Value = '  -3.40%  '

$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '18.91'
$ws.Range("E29").Value = '  -0.95%  '

$ws.Range("E30").Value = '  -4.78%  '

$ws.Range("E31").Value = '  +0.89%  '

$ws.Range("E32").Value = '  -2.97%  '

$ws.Range("E33").Value = '  +4.30%  '

$ws.Range("E34").Value = '  -2.85%  '

$ws.Range("E35").Value = '  -2.44%  '

$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '6.35'
$ws.Range("E36").Value = '  +5.29%  '

$ws.Range("E37").Value = '  -4.17%  '

$ws.Range("E38").Value = '  -3.13%  '

$ws.Range("E39").Value = '  +0.03%  '

$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '1.531.10'
$ws.Range("E40").Value = '  +3.29%  '

$ws.Range("E41").Value = '  -1.07%  '

$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '96.64'
$ws.Range("E42").Value = '  -1.94%  '

$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '16.85'
$ws.Range("E43").Value = '  -0.09%  '

$ws.Range("E44").Value = '  -0.45%  '

$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '0.0915'
$ws.Range("E45").Value = '  -3.43%  '

$ws.Range("B46").Value = 'FTXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '4.18'
$ws.Range("E46").Value = '  +3.18%  '

$ws.Range("B47").Value = 'TrustWalletToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '1.11'
$ws.Range("E47").Value = '  -2.06%  '

$ws.Range("E48").Value = '  -1.85%  '

$ws.Range("B49").Value = 'MXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '2.95'
$ws.Range("E49").Value = '  -0.44%  '

$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '7.12'
$ws.Range("E50").Value = '  +0.30%  '

$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '2.219.29'
$ws.Range("E51").Value = '  -1.59%  '
